$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: Insert a new paragraph right after the document title (the
# Heading1 paragraph "Play Cat Gangster Free Slot Game Review | High 5
# Games") containing the meta description, with "Meta description" in bold
# followed by the plain-text description.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$insertionPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)

# Using InsertXML lets us create the paragraph with the exact run layout
# (a leading empty run, a bold run, then a plain run) that matches the
# rest of the document's authoring style. InsertXML requires the inserted
# content to end on a paragraph boundary, so we include a second, throwaway
# empty paragraph after the one we actually want and remove it afterwards.
$metaXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r/>
<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
<w:r><w:t>: Read our review of Cat Gangster, a fun online slot game by High 5 Games. Play for free and discover the game's multiway payline mechanism and bonus features.</w:t></w:r>
</w:p>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
[void]$insertionPoint.InsertXML($metaXml)

# Remove the spare empty paragraph InsertXML left behind (it is paragraph 3:
# title, meta-description, spare), by deleting its paragraph mark so it
# merges away without touching its neighbours.
$sparePara = $d.Paragraphs.Item(3)
$d.Range($sparePara.Range.Start, $sparePara.Range.End).Delete()

# ---------------------------------------------------------------------------
# Part 2: Near the end of the document, remove the paragraph duplicating the
# bold title text, and rewrite the following italic paragraph's text to hold
# the new image-generation prompt instead of the old meta description text.
# ---------------------------------------------------------------------------
$boldTitleText = "Play Cat Gangster Free Slot Game Review | High 5 Games"
$boldPara = $null
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd() -eq $boldTitleText) {
        $boldPara = $candidate
        break
    }
}
if ($boldPara -ne $null) {
    $d.Range($boldPara.Range.Start, $boldPara.Range.End).Delete()
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
# Exclude the trailing paragraph-mark character from the range so only the
# visible text is replaced, leaving the paragraph's run/format structure
# (leading empty run + italic run) intact.
$italicTextRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)

$newPrompt = "Prompt: As the game `"Cat Gangster`" features a happy Maya warrior with glasses, create a cartoon-style feature image that showcases this character. The image should have a fun and playful vibe and represent the game's gangster theme. The background should feature some seedy downtown alleyways, as this is where the action takes place. The image should also incorporate some of the game's symbols, such as the cat gangster, the logo symbol, and the Kit Kat Club symbol. Overall, the image should capture the essence of the game and entice players to give it a spin."

$italicTextRange.Text = $newPrompt
